# Regenerate orders with updated distance/size labels.
# Distance codes: D80 -> D86, D51 -> D55, D64 -> D69
# Size code:      S30 -> S31  (S20 / S25 unchanged)
#
# These substitutions touch every cell that embeds the old codes:
# the Condition column (e.g. Face06_D80_S25), the Filename_Left /
# Filename_Right columns (e.g. Face06_D80_S25_l.png, Fixation_D80_l.png),
# and the standalone Distance / Size columns (D80, S30, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

[void]$usedRange.Replace("D80", "D86")
[void]$usedRange.Replace("D51", "D55")
[void]$usedRange.Replace("D64", "D69")
[void]$usedRange.Replace("S30", "S31")
